$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update quarter header labels (rolling window: drop oldest quarter, append newest) ---
$quarters = @(
  "فصل سوم منتهی به 1399/09",
  "فصل چهارم منتهی به 1399/12",
  "فصل اول منتهی به 1400/03",
  "فصل دوم منتهی به 1400/06",
  "فصل سوم منتهی به 1400/09",
  "فصل چهارم منتهی به 1400/12",
  "فصل اول منتهی به 1401/03",
  "فصل دوم منتهی به 1401/06",
  "فصل سوم منتهی به 1401/09",
  "فصل چهارم منتهی به 1401/12"
)
$headerRows = @(8, 19, 29, 39, 48, 58)
foreach ($r in $headerRows) {
  for ($i = 0; $i -lt $quarters.Length; $i++) {
    $col = 5 + $i  # column E = 5
    $ws.Cells.Item($r, $col).Value = $quarters[$i]
  }
}

# --- Update data rows (rolling window of numeric/placeholder values) ---
$dataRows = @{
  10 = @(1642324, 1387831, 1717716, 1474261, 1042134, 1240149, 1638788, 1553189, 1347784, 816525)
  11 = @(765006, 647420, 787488, 678231, 483987, 564382, 754136, 680043, 672641, 381232)
  12 = @("-", "-", "-", 0, 0, 0, 0, 0, 0, 0)
  13 = @("-", 0, "-", "-", "-", "-", "-", "-", "-", "-")
  14 = @("-", "-", 0, "-", "-", "-", "-", "-", "-", "-")
  15 = @(2407330, 2035251, 2505204, 2152492, 1526121, 1804531, 2392924, 2233232, 2020425, 1197757)
  21 = @(325882, 226224, 332803, 286947, 158098, 212823, 306950, 246208, 212779, 142146)
  22 = @(0, 0, "-", "-", "-", "-", "-", "-", "-", "-")
  23 = @(926363, 555687, 826763, 714034, 476135, 431025, 783665, 742308, 544049, 390117)
  24 = @(0, 0, "-", 0, 0, 0, 0, 0, 0, 0)
  25 = @(1252245, 781911, 1159566, 1000981, 634233, 643848, 1090615, 988516, 756828, 532263)
  31 = @(219036, 314486, 266649, 229908, 126672, 170518, 553354, 443851, 383587, 256254)
  32 = @(0, 0, "-", "-", "-", "-", "-", "-", "-", "-")
  33 = @(57171847, 39459201, 54418716, 53937022, 45416540, 38063429, 68953860, 57626938, 45973788, 40021974)
  34 = @(0, 0, "-", 0, 0, 0, 0, 0, 0, 0)
  35 = @(57390883, 39773687, 54685365, 54166930, 45543212, 38233947, 69507214, 58070789, 46357375, 40278228)
  41 = @(672133, 1390153, 801222, 801221, 801225, 801220, 1802750, 1802748, 1802748, 1802752)
  43 = @(61716462, 71009761, 65821422, 75537036, 95385846, 88309098, 87988927, 77632133, 84502997, 102589722)
  50 = @(-46873, -22894, -36680, -22125, -38102, -50020, -49050, -56139, -35357, -56624)
  51 = @(0, 0, "-", "-", "-", "-", "-", "-", "-", "-")
  52 = @(-21839693, -20215427, -20856455, -28914647, -32248722, -50057672, -44092487, -45155679, -29998663, -27906611)
  53 = @(0, 0, "-", 0, 0, 0, 0, 0, 0, 0)
  54 = @(-21886566, -20238321, -20893135, -28936772, -32286824, -50107692, -44141537, -45211818, -30034020, -27963235)
  60 = @(172163, 291592, 229969, 207783, 88570, 120498, 504304, 387712, 348230, 199630)
  61 = @(0, 0, "-", "-", "-", "-", "-", "-", "-", "-")
  62 = @(35332154, 19243774, 33562261, 25021365, 13167818, -11994243, 24861373, 12471259, 15975125, 12115363)
  63 = @(0, 0, "-", 0, 0, 0, 0, 0, 0, 0)
  64 = @(35504317, 19535366, 33792230, 25229148, 13256388, -11873745, 25365677, 12858971, 16323355, 12314993)
}

foreach ($r in $dataRows.Keys) {
  $vals = $dataRows[$r]
  for ($i = 0; $i -lt $vals.Length; $i++) {
    $col = 5 + $i
    $ws.Cells.Item([int]$r, $col).Value = $vals[$i]
  }
}

Write-Host "edit applied"